# Update NATMI LR-pair TPM-derived metrics (Il6-Il6ra, YoungD0) with new TPM values.
# The sheet holds one row per Sending/Target cluster pair; columns G-J are the
# ligand-side average/total expression + specificity, M-P the receptor-side
# equivalents, and Q-T the derived edge weights/specificities. K/L (receptor
# detection counts/rates) only change where the TPM re-run flipped the
# detection threshold for the MuSCs receptor-expressing-cell count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2022703333333333
$ws.Range("H2").Value = 0.606811
$ws.Range("I2").Value = 0.01168815774551004
$ws.Range("J2").Value = 0.01168815774551004
$ws.Range("M2").Value = 0.7893693333333335
$ws.Range("N2").Value = 2.368108
$ws.Range("O2").Value = 0.1840020898203156
$ws.Range("P2").Value = 0.1840020898203156
$ws.Range("Q2").Value = 0.1596659981764445
$ws.Range("R2").Value = 1.436993983588
$ws.Range("S2").Value = 0.002150645451323356
$ws.Range("T2").Value = 0.002150645451323355

# Row 3
$ws.Range("G3").Value = 0.2022703333333333
$ws.Range("H3").Value = 0.606811
$ws.Range("I3").Value = 0.01168815774551004
$ws.Range("J3").Value = 0.01168815774551004
$ws.Range("O3").Value = 0.6904048063380857
$ws.Range("P3").Value = 0.6904048063380857
$ws.Range("Q3").Value = 0.5990919595393334
$ws.Range("R3").Value = 5.391827635854
$ws.Range("S3").Value = 0.008069560284737856
$ws.Range("T3").Value = 0.008069560284737855

# Row 4
$ws.Range("G4").Value = 0.2022703333333333
$ws.Range("H4").Value = 0.606811
$ws.Range("I4").Value = 0.01168815774551004
$ws.Range("J4").Value = 0.01168815774551004
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5387946666666666
$ws.Range("N4").Value = 1.616384
$ws.Range("O4").Value = 0.1255931038415988
$ws.Range("P4").Value = 0.1255931038415988
$ws.Range("Q4").Value = 0.1089821768248889
$ws.Range("R4").Value = 0.9808395914239999
$ws.Range("S4").Value = 0.00146795200944883
$ws.Range("T4").Value = 0.00146795200944883

# Row 5
$ws.Range("I5").Value = 0.8358439174604506
$ws.Range("J5").Value = 0.8358439174604506
$ws.Range("M5").Value = 0.7893693333333335
$ws.Range("N5").Value = 2.368108
$ws.Range("O5").Value = 0.1840020898203156
$ws.Range("P5").Value = 0.1840020898203156
$ws.Range("Q5").Value = 11.41804006301156
$ws.Range("R5").Value = 102.762360567104
$ws.Range("S5").Value = 0.1537970275763223
$ws.Range("T5").Value = 0.1537970275763223

# Row 6
$ws.Range("I6").Value = 0.8358439174604506
$ws.Range("J6").Value = 0.8358439174604506
$ws.Range("O6").Value = 0.6904048063380857
$ws.Range("P6").Value = 0.6904048063380857
$ws.Range("S6").Value = 0.5770706579631493
$ws.Range("T6").Value = 0.5770706579631493

# Row 7
$ws.Range("I7").Value = 0.8358439174604506
$ws.Range("J7").Value = 0.8358439174604506
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5387946666666666
$ws.Range("N7").Value = 1.616384
$ws.Range("O7").Value = 0.1255931038415988
$ws.Range("P7").Value = 0.1255931038415988
$ws.Range("Q7").Value = 7.793536979399111
$ws.Range("R7").Value = 70.141832814592
$ws.Range("S7").Value = 0.1049762319209791
$ws.Range("T7").Value = 0.1049762319209791

# Row 8
$ws.Range("G8").Value = 2.638545666666667
$ws.Range("H8").Value = 7.915637
$ws.Range("I8").Value = 0.1524679247940394
$ws.Range("J8").Value = 0.1524679247940394
$ws.Range("M8").Value = 0.7893693333333335
$ws.Range("N8").Value = 2.368108
$ws.Range("O8").Value = 0.1840020898203156
$ws.Range("P8").Value = 0.1840020898203156
$ws.Range("Q8").Value = 2.082787033866223
$ws.Range("R8").Value = 18.745083304796
$ws.Range("S8").Value = 0.02805441679266997
$ws.Range("T8").Value = 0.02805441679266997

# Row 9
$ws.Range("G9").Value = 2.638545666666667
$ws.Range("H9").Value = 7.915637
$ws.Range("I9").Value = 0.1524679247940394
$ws.Range("J9").Value = 0.1524679247940394
$ws.Range("O9").Value = 0.6904048063380857
$ws.Range("P9").Value = 0.6904048063380857
$ws.Range("Q9").Value = 7.814944820268668
$ws.Range("R9").Value = 70.334503382418
$ws.Range("S9").Value = 0.1052645880901986
$ws.Range("T9").Value = 0.1052645880901986

# Row 10
$ws.Range("G10").Value = 2.638545666666667
$ws.Range("H10").Value = 7.915637
$ws.Range("I10").Value = 0.1524679247940394
$ws.Range("J10").Value = 0.1524679247940394
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.5387946666666666
$ws.Range("N10").Value = 1.616384
$ws.Range("O10").Value = 0.1255931038415988
$ws.Range("P10").Value = 0.1255931038415988
$ws.Range("Q10").Value = 1.421634332956444
$ws.Range("R10").Value = 12.794708996608
$ws.Range("S10").Value = 0.01914891991117088
$ws.Range("T10").Value = 0.01914891991117088
